$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '26.207.32'
$ws.Range('E2').Value = '  +3.17%  '

Set-TextValue 'D3' '1.604.03'
$ws.Range('E3').Value = '  +1.92%  '

$ws.Range('E4').Value = '  -0.79%  '

Set-TextValue 'D5' '212.64'
$ws.Range('E5').Value = '  +2.66%  '

$ws.Range('E6').Value = '  -0.75%  '

$ws.Range('E7').Value = '  +0.86%  '

Set-TextValue 'D8' '0.249'
$ws.Range('E8').Value = '  +1.81%  '

$ws.Range('E9').Value = '  +1.72%  '

Set-TextValue 'D10' '18.08'
$ws.Range('E10').Value = '  +2.23%  '

Set-TextValue 'D11' '0.0819'
$ws.Range('E11').Value = '  +4.59%  '

Set-TextValue 'D12' '1.833.18'
$ws.Range('E12').Value = '  +2.15%  '

Set-TextValue 'D13' '1.605.95'
$ws.Range('E13').Value = '  +1.80%  '

$ws.Range('E14').Value = '  -0.52%  '

$ws.Range('E15').Value = '  +0.88%  '

Set-TextValue 'D16' '26.197.87'
$ws.Range('E16').Value = '  +3.07%  '

Set-TextValue 'D17' '60.58'
$ws.Range('E17').Value = '  +0.98%  '

$ws.Range('E18').Value = '  +2.85%  '

$ws.Range('E19').Value = '  -0.68%  '

Set-TextValue 'D20' '197.98'
$ws.Range('E20').Value = '  +6.22%  '

Set-TextValue 'D21' '4.24'
$ws.Range('E21').Value = '  +2.45%  '

Set-TextValue 'D22' '9.39'
$ws.Range('E22').Value = '  +1.14%  '

$ws.Range('E23').Value = '  +1.94%  '

$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D24' '1.78'
$ws.Range('E24').Value = '  +4.56%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D25' '0.129'
$ws.Range('E25').Value = '  +1.37%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '142.45'
$ws.Range('E26').Value = '  +1.04%  '

Set-TextValue 'D27' '1.01'
$ws.Range('E27').Value = '  -0.87%  '

$ws.Range('E28').Value = '  +1.89%  '

$ws.Range('E29').Value = '  -0.01%  '

$ws.Range('E30').Value = '  +0.30%  '

Set-TextValue 'D31' '0.0470'
$ws.Range('E31').Value = '  +1.50%  '

$ws.Range('E32').Value = '  +2.56%  '

Set-TextValue 'D33' '3.00'
$ws.Range('E33').Value = '  +0.34%  '

Set-TextValue 'D34' '1.49'
$ws.Range('E34').Value = '  +2.18%  '

$ws.Range('E35').Value = '  -1.58%  '

Set-TextValue 'D36' '1.106.70'
$ws.Range('E36').Value = '  +2.12%  '

$ws.Range('E37').Value = '  -0.51%  '

$ws.Range('E38').Value = '  +1.67%  '

Set-TextValue 'D39' '2.32'
$ws.Range('E39').Value = '  -0.31%  '

$ws.Range('E40').Value = '  +1.50%  '

Set-TextValue 'D41' '0.500'
$ws.Range('E41').Value = '  +1.55%  '

$ws.Range('E42').Value = '  +6.73%  '

Set-TextValue 'D43' '1.744.79'
$ws.Range('E43').Value = '  +1.83%  '

Set-TextValue 'D44' '5.10'
$ws.Range('E44').Value = '  +0.74%  '

Set-TextValue 'D45' '92.56'
$ws.Range('E45').Value = '  -2.27%  '

Set-TextValue 'D46' '0.0₆0108'
$ws.Range('E46').Value = '  +0.68%  '

$ws.Range('E47').Value = '  +9.10%  '

Set-TextValue 'D48' '53.69'
$ws.Range('E48').Value = '  +1.81%  '

$ws.Range('E49').Value = '  -0.02%  '

$ws.Range('E50').Value = '  +0.33%  '

Set-TextValue 'D51' '1.00'
$ws.Range('E51').Value = '  -0.49%  '
